# Auto-generated edit script: update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, $Addr, $Val)
    $Sheet.Range($Addr).NumberFormat = "@"
    $Sheet.Range($Addr).Value = $Val
    $Sheet.Range($Addr).Style = "Normal"
}

Set-TextCell $ws "D2" "65.534.89"
Set-TextCell $ws "E2" "  -4.60%  "
Set-TextCell $ws "D3" "3.261.18"
Set-TextCell $ws "E3" "  -5.65%  "
Set-TextCell $ws "D4" "0.999"
Set-TextCell $ws "E4" "  -0.02%  "
Set-TextCell $ws "D5" "554.58"
Set-TextCell $ws "E5" "  -3.07%  "
Set-TextCell $ws "D6" "180.10"
Set-TextCell $ws "E6" "  -5.36%  "
Set-TextCell $ws "E7" "  -0.10%  "
Set-TextCell $ws "D8" "0.587"
Set-TextCell $ws "E8" "  -2.87%  "
Set-TextCell $ws "D9" "3.262.44"
Set-TextCell $ws "E9" "  -5.32%  "
Set-TextCell $ws "D10" "0.185"
Set-TextCell $ws "E10" "  -8.03%  "
Set-TextCell $ws "D11" "0.585"
Set-TextCell $ws "E11" "  -4.57%  "
Set-TextCell $ws "D12" "47.32"
Set-TextCell $ws "E12" "  -7.16%  "
Set-TextCell $ws "D13" "0.0000265"
Set-TextCell $ws "E13" "  -6.47%  "
Set-TextCell $ws "B14" "BitcoinCash"
Set-TextCell $ws "C14" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell $ws "D14" "633.77"
Set-TextCell $ws "E14" "  -0.30%  "
Set-TextCell $ws "B15" "Polkadot"
Set-TextCell $ws "C15" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws "D15" "8.55"
Set-TextCell $ws "E15" "  -5.48%  "
Set-TextCell $ws "D16" "3.773.72"
Set-TextCell $ws "D17" "65.340.96"
Set-TextCell $ws "E17" "  -4.65%  "
Set-TextCell $ws "B18" "Chainlink"
Set-TextCell $ws "C18" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell $ws "D18" "17.74"
Set-TextCell $ws "E18" "  -1.83%  "
Set-TextCell $ws "B19" "TRON"
Set-TextCell $ws "C19" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell $ws "D19" "0.116"
Set-TextCell $ws "E19" "  -3.18%  "
Set-TextCell $ws "D20" "3.249.41"
Set-TextCell $ws "E20" "  -6.09%  "
Set-TextCell $ws "D21" "11.35"
Set-TextCell $ws "E21" "  -7.45%  "
Set-TextCell $ws "D22" "0.901"
Set-TextCell $ws "E22" "  -3.76%  "
Set-TextCell $ws "D23" "17.68"
Set-TextCell $ws "E23" "  -0.88%  "
Set-TextCell $ws "D24" "106.09"
Set-TextCell $ws "E24" "  +7.07%  "
Set-TextCell $ws "D25" "4.97"
Set-TextCell $ws "E25" "  -6.39%  "
Set-TextCell $ws "E26" "  -7.04%  "
Set-TextCell $ws "D27" "2.67"
Set-TextCell $ws "E27" "  -5.60%  "
Set-TextCell $ws "D28" "9.51"
Set-TextCell $ws "E28" "  -2.28%  "
Set-TextCell $ws "D29" "8.71"
Set-TextCell $ws "E29" "  -5.04%  "
Set-TextCell $ws "D30" "30.25"
Set-TextCell $ws "E30" "  -6.12%  "
Set-TextCell $ws "D31" "4.07"
Set-TextCell $ws "E31" "  -2.19%  "
Set-TextCell $ws "D32" "6.32"
Set-TextCell $ws "E32" "  -5.58%  "
Set-TextCell $ws "D33" "11.03"
Set-TextCell $ws "E33" "  -4.29%  "
Set-TextCell $ws "D34" "549.76"
Set-TextCell $ws "E34" "  +9.82%  "
Set-TextCell $ws "E35" "  -2.56%  "
Set-TextCell $ws "B36" "OKB"
Set-TextCell $ws "C36" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws "D36" "57.03"
Set-TextCell $ws "E36" "  -6.31%  "
Set-TextCell $ws "B37" "Dai"
Set-TextCell $ws "C37" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell $ws "D37" "1.00"
Set-TextCell $ws "E37" "  +0.06%  "
Set-TextCell $ws "D38" "3.621.95"
Set-TextCell $ws "E38" "  -1.05%  "
Set-TextCell $ws "D39" "3.72"
Set-TextCell $ws "E39" "  +8.92%  "
Set-TextCell $ws "D40" "3.41"
Set-TextCell $ws "E40" "  -2.32%  "
Set-TextCell $ws "D41" "2.74"
Set-TextCell $ws "E41" "  -4.26%  "
Set-TextCell $ws "D42" "0.130"
Set-TextCell $ws "E42" "  -1.68%  "
Set-TextCell $ws "D43" "0.0₃0712"
Set-TextCell $ws "E43" "  -8.05%  "
Set-TextCell $ws "D44" "31.89"
Set-TextCell $ws "E44" "  -6.78%  "
Set-TextCell $ws "D45" "0.338"
Set-TextCell $ws "E45" "  -7.55%  "
Set-TextCell $ws "D46" "3.30"
Set-TextCell $ws "E46" "  -1.31%  "
Set-TextCell $ws "D47" "0.0413"
Set-TextCell $ws "E47" "  -4.85%  "
Set-TextCell $ws "D48" "2.61"
Set-TextCell $ws "E48" "  -6.62%  "
Set-TextCell $ws "D49" "0.129"
Set-TextCell $ws "E49" "  -3.34%  "
Set-TextCell $ws "D50" "0.998"
Set-TextCell $ws "E50" "  -0.15%  "
Set-TextCell $ws "E51" "  +2.01%  "

Write-Output "Applied 107 cell updates"
